$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.78125,
    0.71875,
    0.53125,
    0.46875,
    0.421875,
    0.421875,
    0.390625,
    0.40625,
    0.3125,
    0.390625,
    0.390625,
    0.390625,
    0.359375,
    0.375,
    0.40625,
    0.40625,
    0.390625,
    0.40625,
    0.40625,
    0.4375,
    0.375,
    0.375,
    0.515625,
    0.4375,
    0.40625,
    0.34375,
    0.40625,
    0.375,
    0.4375,
    0.390625,
    0.390625,
    0.453125,
    0.421875,
    0.421875,
    0.421875,
    0.453125,
    0.546875,
    0.4375,
    0.375,
    0.375,
    0.390625,
    0.53125,
    0.421875,
    0.359375,
    0.375,
    0.390625,
    0.40625,
    0.40625,
    0.40625,
    0.40625,
    0.40625,
    0.40625,
    0.40625,
    0.40625,
    0.40625,
    0.421875,
    0.421875,
    0.421875,
    0.421875,
    0.4375,
    0.4375,
    0.4375,
    0.453125,
    0.421875,
    0.390625,
    0.390625,
    0.375,
    0.390625,
    0.359375,
    0.375,
    0.375,
    0.375,
    0.375,
    0.390625,
    0.40625,
    0.390625,
    0.390625,
    0.390625,
    0.390625,
    0.390625,
    0.40625,
    0.390625,
    0.375,
    0.375,
    0.375,
    0.390625,
    0.390625,
    0.390625,
    0.390625,
    0.390625,
    0.390625,
    0.390625,
    0.390625,
    0.40625,
    0.40625,
    0.40625,
    0.40625,
    0.40625,
    0.40625,
    0.40625,
    0.40625,
    0.390625,
    0.25,
    0.328125,
    0.25,
    0.203125,
    0.28125,
    0.28125,
    0.359375,
    0.421875,
    0.328125,
    0.28125,
    0.21875,
    0.296875,
    0.3125,
    0.28125,
    0.3114754098360656
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$newLabel = "<__main__.DisplayOutputs object at 0x7f53f262f3d0>"
for ($row = 102; $row -le 118; $row++) {
    $ws.Cells.Item($row, 1).Value = $newLabel
}

Write-Host "Updated $($values.Length) cells in column B and relabeled A102:A118"